# Added some logs into Respin
# Update the per-row reel data on Sheet1 (rows 2-22, columns A-F) to reflect
# the new log values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @(901, 16, 15, 45, 60, 60)
    3  = @(101, 9, 30, 15, 60, 15)
    4  = @(301, 6, 45, 30, 60, 45)
    5  = @(701, 3, 90, 45, 97, 15)
    6  = @(801, 3, 67, 65, 52, 45)
    7  = @(1203, 3, 15, 15, 15, 15)
    8  = @(1001, 18, 30, 75, 60, 72)
    9  = @(401, 9, 48, 67, 75, 45)
    10 = @(1201, 2, 10, 10, 10, 10)
    11 = @(1202, 2, 10, 10, 10, 10)
    12 = @(902, 1, 0, 0, 0, 0)
    13 = @(501, 9, 52, 30, 75, 45)
    14 = @(601, 9, 60, 67, 60, 42)
    15 = @(201, 9, 30, 15, 45, 30)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(1, 0, 2, 2, 2, 2)
    18 = @(2, 0, 2, 2, 2, 2)
    19 = @(1101, 0, 15, 30, 30, 0)
    20 = @(802, 0, 4, 5, 4, 0)
    21 = @(3, 0, 3, 3, 3, 3)
    22 = @(602, 0, 0, 4, 0, 9)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
